$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-arrange row 4 header values:
# Before: A4=name, B4=amount, C4=price,       D4=mass
# After:  A4=name, B4=price,  C4=vendor_code, D4=mass, E4=amount
$ws.Range("B4").Value = "price"
$ws.Range("C4").Value = "vendor_code"
$ws.Range("D4").Value = "mass"
$ws.Range("E4").Value = "amount"

# Adjust column widths for C and D to match new content widths
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 9

# Update the active selection to A4 (matches saved view state)
$ws.Range("A4").Select()
